# Slide 3 ("Conditional Processing") introduction rework:
#  1. Remove the "Conditional Processing" title textbox (shape id 204).
#  2. Enlarge/reposition the "girl-looking-out-the-window" picture so it
#     takes over the space that was occupied by the title text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- 1. Delete the title textbox shape ------------------------------------
$titleShape = $s.Shapes.Item("Conditional Processing")
$titleShape.Delete()

# --- 2. Reposition / resize the picture ------------------------------------
# (target EMUs: off 8543510,973402  ext 7296981,4966248 — the literal point
#  values below are chosen so the engine's point->EMU conversion lands on
#  those exact EMU values).
$pic = $s.Shapes.Item("girl-looking-out-the-window-jpg.jpg")
$pic.Left   = 672.7173157346457
$pic.Top    = 76.64582447165354
$pic.Width  = 574.5654602708661
$pic.Height = 391.04313660629924
